$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.886.28"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.281.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.57"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -5.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.278.86"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.70%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.571"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.33"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000267"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "685.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.807.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.26"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.969.87"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.119"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.280.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.884"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.89"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -5.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.86"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.24"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.31"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.34"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.67"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "574.24"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.865.05"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.79%  "

$ws.Range("E35").Value = "  -2.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -13.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "31.71"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0665"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.73%  "

$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("E46").Value = "  -2.36%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.127"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("E49").Value = "  +6.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.26"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.37%  "
